$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column L into column M for rows 2 through 33
for ($r = 2; $r -le 33; $r++) {
    $ws.Range("L$r").Copy() | Out-Null
    $ws.Range("M$r").PasteSpecial(-4122) | Out-Null
}

# Set the new values for column M (year 2021 data)
$ws.Range("M3").Value = 2021
$ws.Range("M4").Value = 16.60175395812114
$ws.Range("M5").Value = 7.3576531459900787
$ws.Range("M6").Value = 25.982831460241147
$ws.Range("M7").Value = 4.6978046797362003
$ws.Range("M8").Value = 0.36820478077087354
$ws.Range("M9").Value = 8.8695886639561206
$ws.Range("M10").Value = 4.8763798385289059
$ws.Range("M11").Value = 0.63362020488109327
$ws.Range("M12").Value = 9.0604897038469581
$ws.Range("M13").Value = 6.1518067459522099
$ws.Range("M14").Value = 2.7662407973096332
$ws.Range("M15").Value = 9.5668603955896767
$ws.Range("M16").Value = 19.43893106341184
$ws.Range("M17").Value = 6.948932296552635
$ws.Range("M18").Value = 31.476235442241109
$ws.Range("M19").Value = 17.246785826277829
$ws.Range("M20").Value = 1.7299267231872171
$ws.Range("M21").Value = 32.417697807858893
$ws.Range("M22").Value = 13.57738374823861
$ws.Range("M23").Value = 1.4831184047578438
$ws.Range("M24").Value = 25.424775353949197
$ws.Range("M25").Value = 28.259337058396849
$ws.Range("M26").Value = 18.681751343880823
$ws.Range("M27").Value = 38.14301557851968
$ws.Range("M28").Value = 29.550909035673744
$ws.Range("M29").Value = 19.441305483663228
$ws.Range("M30").Value = 40.966897178968502
$ws.Range("M31").Value = 17.68853538926977
$ws.Range("M32").Value = 4.1574114768313395
$ws.Range("M33").Value = 31.970511904314137

# Update the active selection to match the authored workbook
$ws.Range("N4").Select() | Out-Null

